$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Pdpn"
$ws.Range("C2").Value = "Clec1b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 53.570086
$ws.Range("H2").Value = 160.710258
$ws.Range("I2").Value = 0.927800950569834
$ws.Range("J2").Value = 0.927800950569834
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.450881
$ws.Range("N2").Value = 4.352643
$ws.Range("O2").Value = 0.3369514394030648
$ws.Range("P2").Value = 0.3369514394030648
$ws.Range("Q2").Value = 77.72381994576601
$ws.Range("R2").Value = 699.5143795118942
$ws.Range("S2").Value = 0.3126238657740373
$ws.Range("T2").Value = 0.3126238657740373

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Pdpn"
$ws.Range("C3").Value = "Clec1b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 53.570086
$ws.Range("H3").Value = 160.710258
$ws.Range("I3").Value = 0.927800950569834
$ws.Range("J3").Value = 0.927800950569834
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.609510666666667
$ws.Range("N3").Value = 7.828531999999999
$ws.Range("O3").Value = 0.6060306636250556
$ws.Range("P3").Value = 0.6060306636250556
$ws.Range("Q3").Value = 139.7917108312507
$ws.Range("R3").Value = 1258.125397481256
$ws.Range("S3").Value = 0.5622758257857939
$ws.Range("T3").Value = 0.5622758257857939

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Pdpn"
$ws.Range("C4").Value = "Clec1b"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 53.570086
$ws.Range("H4").Value = 160.710258
$ws.Range("I4").Value = 0.927800950569834
$ws.Range("J4").Value = 0.927800950569834
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2455136666666667
$ws.Range("N4").Value = 0.736541
$ws.Range("O4").Value = 0.05701789697187955
$ws.Range("P4").Value = 0.05701789697187955
$ws.Range("Q4").Value = 13.15218823750867
$ws.Range("R4").Value = 118.369694137578
$ws.Range("S4").Value = 0.05290125901000271
$ws.Range("T4").Value = 0.05290125901000271

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Pdpn"
$ws.Range("C5").Value = "Clec1b"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.168684333333334
$ws.Range("H5").Value = 12.506053
$ws.Range("I5").Value = 0.07219904943016597
$ws.Range("J5").Value = 0.07219904943016596
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.450881
$ws.Range("N5").Value = 4.352643
$ws.Range("O5").Value = 0.3369514394030648
$ws.Range("P5").Value = 0.3369514394030648
$ws.Range("Q5").Value = 6.048264894231002
$ws.Range("R5").Value = 54.43438404807901
$ws.Range("S5").Value = 0.02432757362902745
$ws.Range("T5").Value = 0.02432757362902744

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Pdpn"
$ws.Range("C6").Value = "Clec1b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.168684333333334
$ws.Range("H6").Value = 12.506053
$ws.Range("I6").Value = 0.07219904943016597
$ws.Range("J6").Value = 0.07219904943016596
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.609510666666667
$ws.Range("N6").Value = 7.828531999999999
$ws.Range("O6").Value = 0.6060306636250556
$ws.Range("P6").Value = 0.6060306636250556
$ws.Range("Q6").Value = 10.87822623379956
$ws.Range("R6").Value = 97.904036104196
$ws.Range("S6").Value = 0.04375483783926168
$ws.Range("T6").Value = 0.04375483783926167

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Pdpn"
$ws.Range("C7").Value = "Clec1b"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.168684333333334
$ws.Range("H7").Value = 12.506053
$ws.Range("I7").Value = 0.07219904943016597
$ws.Range("J7").Value = 0.07219904943016596
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2455136666666667
$ws.Range("N7").Value = 0.736541
$ws.Range("O7").Value = 0.05701789697187955
$ws.Range("P7").Value = 0.05701789697187955
$ws.Range("Q7").Value = 1.023468975852556
$ws.Range("R7").Value = 9.211220782673001
$ws.Range("S7").Value = 0.004116637961876842
$ws.Range("T7").Value = 0.004116637961876841

$ws.Range("A8:T10").EntireRow.Delete() | Out-Null

Write-Output "done"
